# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row -> new value for sheet "展览"
$exhibitUpdates = @{
    3  = 269
    4  = 76
    5  = 263
    6  = 267
    7  = 91
    12 = 111
    13 = 2390
    17 = 553
    21 = 50
    22 = 1906
    23 = 4042
    26 = 1186
    27 = 229
    28 = 2094
    30 = 468
    32 = 114
    33 = 289
    35 = 461
    36 = 696
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

# Row -> new value for sheet "全部类型"
$allUpdates = @{
    3  = 269
    4  = 76
    5  = 263
    6  = 267
    7  = 91
    12 = 111
    13 = 2390
    18 = 553
    22 = 50
    23 = 1906
    24 = 4042
    27 = 1186
    28 = 229
    29 = 2094
    31 = 468
    33 = 114
    34 = 289
    36 = 461
    37 = 696
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
